# Applies the "branches_1.xlsx" revision: refreshed per-unit R/X values on the
# MT and BT branch sheets, a couple of selection/active-tab changes left
# behind by the author's last interactive session, and clears the stray
# yellow highlight that had been left on MT!C12:D12.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# MT sheet: bus A2 renumbered, and R (col C) / X (col D) recomputed to
# per-unit values for rows 2-13.
# ---------------------------------------------------------------------
$mt = $wb.Worksheets.Item("MT")

$mt.Range("A2").Value = 632

$mtData = @{
    2  = @(0.043157, 0.130493)
    3  = @(0.014393, 0.043519)
    4  = @(0.028764, 0.086973)
    5  = @(0.021578, 0.065246)
    6  = @(0.023608, 0.037783)
    7  = @(0.062813, 0.064015)
    8  = @(0.037688, 0.038409)
    9  = @(0.024885, 0.013399)
    10 = @(0.037688, 0.038409)
    11 = @(0.075523, 0.076563)
    12 = @(0.20341, 0.077637)
    13 = @(0.127131, 0.048523)
}
foreach ($row in $mtData.Keys) {
    $vals = $mtData[$row]
    $mt.Cells.Item($row, 3).Value = $vals[0]
    $mt.Cells.Item($row, 4).Value = $vals[1]
}

# That yellow call-out fill on C12:D12 is gone in the new revision.
$mt.Range("C12:D12").ClearFormats()

# ---------------------------------------------------------------------
# BT sheet: R (col C) / X (col D) recomputed to per-unit values for the
# 556.5 KCM rows (2-21); the 4/0 AWG rows (22+) are unchanged.
# ---------------------------------------------------------------------
$bt = $wb.Worksheets.Item("BT")
for ($row = 2; $row -le 21; $row++) {
    $bt.Cells.Item($row, 3).Value = 0.0019
    $bt.Cells.Item($row, 4).Value = 0.000598
}

# ---------------------------------------------------------------------
# Leftover UI state from the author's last save: MT/BT selections moved,
# and the active tab moved from Trafos to Reg (with a new selection).
# ---------------------------------------------------------------------
$mt.Activate()
$mt.Range("E16").Select()

$bt.Activate()
$bt.Range("D22").Select()

$reg = $wb.Worksheets.Item("Reg")
$reg.Activate()
$reg.Range("C2:C4").Select()
